$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new "2022-Q1" sheet right after "2021-Q4" (mirrors the other
#    per-quarter fund-holding sheets) and before the "总计" summary sheet.
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $q4)
$newSheet.Name = "2022-Q1"

# Bring over the header-row formatting (bold/border/centered) and the index
# cell formatting from the sibling "2021-Q4" sheet so the new sheet matches
# the look of the other quarterly sheets.
$q4.Range("B1:H1").Copy()
$newSheet.Range("B1").PasteSpecial(-4122)

$q4.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$newSheet.Range("A2").Value = 0

# Keep these text-like numeric strings as TEXT (not auto-converted numbers),
# matching the source data's inlineStr cells.
$newSheet.Range("B2:G2").NumberFormat = "@"
$newSheet.Range("B2").Value = "002849"
$newSheet.Range("C2").Value = "金信智能中国2025灵活配置混合"
$newSheet.Range("D2").Value = "1.44"
$newSheet.Range("E2").Value = "82.44"
$newSheet.Range("F2").Value = "4.88"
$newSheet.Range("G2").Value = "0.0703"
$newSheet.Range("H2").Value = 6

# Re-apply the sibling sheet's (unstyled / General) formatting to the data
# row so the text-coercion trick above doesn't leave a stray number format.
$q4.Range("B2:G2").Copy()
$newSheet.Range("B2").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Add a new top row to "总计" summarizing the new 2022-Q1 quarter, pushing
#    the existing rows down by one.
# ---------------------------------------------------------------------------
$tot = $wb.Worksheets.Item("总计")
$tot.Range("A2").EntireRow.Insert()
$tot.Range("A2:D2").ClearFormats()

# Re-apply the same formatting the other index cells (A3:A6) carry.
$tot.Range("A3").Copy()
$tot.Range("A2").PasteSpecial(-4122)

$tot.Range("A2").Value = 0
$tot.Range("A3").Value = 1
$tot.Range("A4").Value = 2
$tot.Range("A5").Value = 3
$tot.Range("A6").Value = 4

$tot.Range("B2").Value = "2022-Q1"
$tot.Range("C2").Value = 1
$tot.Range("D2").Value = 0.07000000000000001
